$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("V4")

# --- 8-2 timing split additions (row 94/95, columns I/J) ---
$ws.Range("I94").Value = 17
$ws.Range("J94").Value = 8

$ws.Range("I95").Formula = "=SUM(I93:I94)"
$ws.Range("J95").Formula = "=SUM(J93:J94)"

# --- corrected B106 value (and its dependent D106 recomputes automatically) ---
$ws.Range("B106").Value = 31906

# --- newly logged splits, rows 107-111 ---
$ws.Range("A107").Value = "Black screen (water)"
$ws.Range("B107").Value = 32382
$ws.Range("C107").Value = 37736
$ws.Range("D107").Formula = "=IF(B107>0,C107-B107,0)"

$ws.Range("A108").Value = "Get flag"
$ws.Range("B108").Value = 32646
$ws.Range("C108").Value = 38007
$ws.Range("D108").Formula = "=IF(B108>0,C108-B108,0)"

$ws.Range("A109").Value = "End level"
$ws.Range("B109").Value = 33164
$ws.Range("C109").Value = 38525
$ws.Range("D109").Formula = "=IF(B109>0,C109-B109,0)"

$ws.Range("A110").Value = "Enter 8-F1"
$ws.Range("B110").Value = 33599
$ws.Range("C110").Value = 39462
$ws.Range("D110").Formula = "=IF(B110>0,C110-B110,0)"

$ws.Range("A111").Value = "1st Move"
$ws.Range("B111").Value = 33826
$ws.Range("C111").Value = 39712
$ws.Range("D111").Formula = "=IF(B111>0,C111-B111,0)"

# --- move the active selection to reflect where editing left off ---
$ws.Range("B112").Select()
